$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 131046824
$ws.Range("Q2").Value = 401653
$ws.Range("R2").Value = 6818054
$ws.Range("Z2").Value = "14:50"
$ws.Range("AB2").Value = "14:50"

# Row 3
$ws.Range("A3").Value = 131046825
$ws.Range("Q3").Value = 401650
$ws.Range("R3").Value = 6818017
$ws.Range("Z3").Value = "14:52"
$ws.Range("AB3").Value = "14:52"

# Row 10
$ws.Range("A10").Value = 131046823
$ws.Range("B10").Value = 79244
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 401661
$ws.Range("R10").Value = 6818064
$ws.Range("Z10").Value = "14:50"
$ws.Range("AB10").Value = "14:50"
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 131046773
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("M11").Value = "äldre spår"
$ws.Range("Q11").Value = 401346
$ws.Range("R11").Value = 6818162
$ws.Range("Z11").Value = "15:23"
$ws.Range("AB11").Value = "15:23"
$ws.Range("AC11").Value = "Äldre ringhack (gran)"

# Row 16
$ws.Range("A16").Value = 131046724
$ws.Range("B16").Value = 79276
$ws.Range("E16").Value = 185
$ws.Range("F16").Value = "Violettgrå tagellav"
$ws.Range("G16").Value = "Bryoria nadvornikiana"
$ws.Range("H16").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q16").Value = 401635
$ws.Range("R16").Value = 6817874
$ws.Range("Z16").Value = "14:58"
$ws.Range("AB16").Value = "14:58"

# Row 17
$ws.Range("A17").Value = 131046708
$ws.Range("B17").Value = 83224
$ws.Range("E17").Value = 6440
$ws.Range("F17").Value = "Vitgrynig nållav"
$ws.Range("G17").Value = "Chaenotheca subroscida"
$ws.Range("H17").Value = "(Eitner) Zahlbr."
$ws.Range("Q17").Value = 401645
$ws.Range("R17").Value = 6818016
$ws.Range("Z17").Value = "14:52"
$ws.Range("AB17").Value = "14:52"

# Row 20
$ws.Range("A20").Value = 131046766
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = "äldre spår"
$ws.Range("N20").Value = ""
$ws.Range("Q20").Value = 401318
$ws.Range("R20").Value = 6818379
$ws.Range("Z20").Value = "15:29"
$ws.Range("AB20").Value = "15:29"
$ws.Range("AC20").Value = "Äldre ringhack (tall)"

# Row 22
$ws.Range("A22").Value = 131046799
$ws.Range("B22").Value = 78256
$ws.Range("E22").Value = 228579
$ws.Range("F22").Value = "Liten svartspik"
$ws.Range("G22").Value = "Chaenothecopsis nana"
$ws.Range("H22").Value = "Tibell"
$ws.Range("K22").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("Q22").Value = 401649
$ws.Range("R22").Value = 6818014
$ws.Range("Z22").Value = "14:52"
$ws.Range("AB22").Value = "14:52"
$ws.Range("AC22").ClearContents()

# Row 25
$ws.Range("A25").Value = 131047014
$ws.Range("B25").Value = 57884
$ws.Range("E25").Value = 100109
$ws.Range("F25").Value = "Tretåig hackspett"
$ws.Range("G25").Value = "Picoides tridactylus"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("M25").Value = "färska spår"
$ws.Range("Q25").Value = 401378
$ws.Range("R25").Value = 6818082
$ws.Range("Z25").Value = "15:21"
$ws.Range("AB25").Value = "15:21"
$ws.Range("AC25").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE25").Value = $true

# Row 26
$ws.Range("A26").Value = 131046832
$ws.Range("Q26").Value = 401350
$ws.Range("R26").Value = 6818162
$ws.Range("Z26").Value = "15:24"
$ws.Range("AB26").Value = "15:24"

# Row 27
$ws.Range("A27").Value = 131046826
$ws.Range("B27").Value = 79244
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("M27").ClearContents()
$ws.Range("Q27").Value = 401661
$ws.Range("R27").Value = 6818064
$ws.Range("Z27").Value = "14:54"
$ws.Range("AB27").Value = "14:54"
$ws.Range("AC27").ClearContents()
$ws.Range("AE27").Value = $false
